$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'310.60"
$ws.Range("E2").Formula = "'1.71%"

$ws.Range("D3").Formula = "'35.58"
$ws.Range("E3").Formula = "'-2.13%"

$ws.Range("D4").Formula = "'5.100"
$ws.Range("E4").Formula = "'1.38%"

$ws.Range("E5").Formula = "'3.45%"

$ws.Range("D6").Formula = "'2.056"
$ws.Range("E6").Formula = "'-3.32%"

$ws.Range("D7").Formula = "'7.949"
$ws.Range("E7").Formula = "'-0.33%"

$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Formula = "'4.126"
$ws.Range("E8").Formula = "'-0.42%"

$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D9").Formula = "'2.974"
$ws.Range("E9").Formula = "'11.76%"

$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D10").Formula = "'0.9248"
$ws.Range("E10").Formula = "'-0.11%"

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Formula = "'0.1130"
$ws.Range("E11").Formula = "'15.74%"

$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D12").Formula = "'0.1919"
$ws.Range("E12").Formula = "'2.64%"

$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").Formula = "'0.09226"
$ws.Range("E13").Formula = "'1.86%"

$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").Formula = "'0.03679"
$ws.Range("E14").Formula = "'2.12%"

$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Formula = "'0.09912"
$ws.Range("E15").Formula = "'-0.05%"

$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Formula = "'0.001440"
$ws.Range("E16").Formula = "'0.37%"

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Formula = "'0.005832"
$ws.Range("E17").Formula = "'3.96%"

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Formula = "'3.481"
$ws.Range("E18").Formula = "'-0.02%"

$ws.Range("D19").Formula = "'0.3402"
$ws.Range("E19").Formula = "'-0.62%"

$ws.Range("E20").Formula = "'-1.37%"

$ws.Range("D21").Formula = "'5.091"
$ws.Range("E21").Formula = "'0.23%"

$ws.Range("D22").Formula = "'0.2204"
$ws.Range("E22").Formula = "'-1.87%"

$ws.Range("D23").Formula = "'0.04535"
$ws.Range("E23").Formula = "'-1.07%"

$ws.Range("D24").Formula = "'0.001225"
$ws.Range("E24").Formula = "'-0.79%"

$ws.Range("D25").Formula = "'0.004812"
$ws.Range("E25").Formula = "'-0.17%"

$ws.Range("D26").Formula = "'0.0001250"
$ws.Range("E26").Formula = "'-3.93%"

$ws.Range("D27").Formula = "'0.0004445"
$ws.Range("E27").Formula = "'-6.18%"

$ws.Range("D39").Formula = "'0.01981"
$ws.Range("E39").Formula = "'3.92%"

$ws.Range("D40").Formula = "'0.04884"
$ws.Range("E40").Formula = "'-0.31%"

$ws.Range("D41").Formula = "'0.007640"
$ws.Range("E41").Formula = "'-2.35%"

$ws.Range("D42").Formula = "'0.009455"
$ws.Range("E42").Formula = "'22.02%"

$ws.Range("D43").Formula = "'0.1387"
$ws.Range("E43").Formula = "'-0.80%"

$ws.Range("D44").Formula = "'0.002121"
$ws.Range("E44").Formula = "'-6.20%"

$ws.Range("D45").Formula = "'0.01163"
$ws.Range("E45").Formula = "'3.55%"

$ws.Range("D46").Formula = "'0.00006555"
$ws.Range("E46").Formula = "'2.17%"

$ws.Range("E47").Formula = "'0.05%"

$ws.Range("D48").Formula = "'180.20"
$ws.Range("E48").Formula = "'247.91%"

$ws.Range("E49").Formula = "'-21.08%"

$ws.Range("E50").Formula = "'0.05%"

$ws.Range("E51").Formula = "'0.05%"
